# ML model retrained with all data
# - J1:K1 used to hold shared-string labels "r"/"s"; now they hold the
#   numeric weights used for every other row (0.3 / 0.6).
# - J2:K51 flip from the old split (0.5 / 0.3) to the new one (0.3 / 0.6).
# - Window/view state (zoom + selection) is updated to match where the
#   author left the workbook after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: was text labels "r"/"s" via shared strings, now numeric.
$ws.Range("J1").Value = 0.3
$ws.Range("K1").Value = 0.6

# All data rows: J goes 0.5 -> 0.3, K goes 0.3 -> 0.6.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.3
    $ws.Cells.Item($r, 11).Value = 0.6
}

# View state: zoom to 90% and leave the K column selected, matching the
# author's saved window position.
$win = $excel.ActiveWindow
$win.Zoom = 90
$win.ScrollRow = 39
$win.ScrollColumn = 1

$ws.Range("K1:K51").Select()
